$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("H6").Value = 0.13772999999999999
$ws.Range("I6").Value = 0.23014999999999999
$ws.Range("J6").Value = 0.15997
